$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TSPARMCD values to drop the trailing "_im" suffix
# (data.table::cast from long to wide cannot have underbars in TSPARMCD values)
$ws.Range("C4").Value = "siteid"
$ws.Range("C5").Value = "adminact"
$ws.Range("C6").Value = "arm"
$ws.Range("C7").Value = "arm"
$ws.Range("C8").Value = "arm"
$ws.Range("C9").Value = "arm"
$ws.Range("C10").Value = "epoch"

# Add a NOTES row explaining the rule
$ws.Range("O2").Value = "DO NOT USE UNDERBARS in any TSPARMCD values. Due to later CAST function from long to wide."
